$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '60.689.51'
$ws.Range('E2').Value = '  -0.43%  '
Set-TextValue 'D3' '2.397.32'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  +0.75%  '
Set-TextValue 'D5' '560.99'
$ws.Range('E5').Value = '  -1.65%  '
Set-TextValue 'D6' '141.24'
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('E7').Value = '  -0.25%  '
Set-TextValue 'D8' '0.540'
$ws.Range('E8').Value = '  +2.60%  '
Set-TextValue 'D9' '2.402.30'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('E13').Value = '  +2.01%  '
Set-TextValue 'D14' '26.15'
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D15' '0.0000168'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D16' '2.812.99'
$ws.Range('E16').Value = '  -0.61%  '
Set-TextValue 'D17' '60.371.83'
$ws.Range('E17').Value = '  -0.82%  '
Set-TextValue 'D18' '2.404.07'
$ws.Range('E18').Value = '  +0.24%  '
Set-TextValue 'D19' '8.15'
$ws.Range('E19').Value = '  +7.26%  '
Set-TextValue 'D20' '10.66'
$ws.Range('E20').Value = '  +0.12%  '
Set-TextValue 'D21' '324.09'
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  -1.10%  '
Set-TextValue 'D26' '64.65'
$ws.Range('E26').Value = '  -0.25%  '
Set-TextValue 'D27' '571.15'
$ws.Range('E27').Value = '  -2.32%  '
Set-TextValue 'D28' '8.09'
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('E29').Value = '  +0.14%  '
Set-TextValue 'D30' '0.0₃0938'
$ws.Range('E30').Value = '  +0.32%  '
Set-TextValue 'D31' '8.05'
$ws.Range('E31').Value = '  +1.85%  '
Set-TextValue 'D32' '1.33'
$ws.Range('E32').Value = '  -1.30%  '
Set-TextValue 'D33' '1.80'
$ws.Range('E33').Value = '  -2.42%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  -0.54%  '
Set-TextValue 'D36' '1.45'
$ws.Range('E36').Value = '  +3.18%  '
Set-TextValue 'D37' '152.05'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('E39').Value = '  -0.50%  '
Set-TextValue 'D40' '18.29'
$ws.Range('E40').Value = '  -0.03%  '
Set-TextValue 'D41' '5.15'
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D42' '0.999'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D43' '2.53'
$ws.Range('E43').Value = '  +6.88%  '
$ws.Range('E44').Value = '  +0.21%  '
Set-TextValue 'D45' '41.65'
$ws.Range('E45').Value = '  +1.12%  '
Set-TextValue 'D46' '0.0₆0278'
$ws.Range('E46').Value = '  -3.44%  '
Set-TextValue 'D47' '141.74'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('E48').Value = '  -0.01%  '
Set-TextValue 'D49' '0.588'
$ws.Range('E49').Value = '  -0.05%  '
Set-TextValue 'D50' '0.0508'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('E51').Value = '  -0.72%  '
